$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K to E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy formatting from column E (old D, now shifted) into the new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Populate new column D with the latest period data
$ws.Range("D7").Value() = 43463
$ws.Range("D8").Value() = 2239200
$ws.Range("D9").Value() = 1317900
$ws.Range("D10").Value() = 921300
$ws.Range("D12").Value() = "NA"
$ws.Range("D13").Value() = 0
$ws.Range("D14").Value() = 600
$ws.Range("D15").Value() = 0
$ws.Range("D17").Value() = 1987900
$ws.Range("D18").Value() = 251300
$ws.Range("D20").Value() = -23900
$ws.Range("D21").Value() = 258900
$ws.Range("D22").Value() = "NA"
$ws.Range("D23").Value() = 227400
$ws.Range("D24").Value() = 27200
$ws.Range("D25").Value() = 0
$ws.Range("D26").Value() = 200200
$ws.Range("D27").Value() = 192500
$ws.Range("D28").Value() = 0
$ws.Range("D29").Value() = 100
$ws.Range("D30").Value() = 0
$ws.Range("D31").Value() = 0
$ws.Range("D32").Value() = 23900
$ws.Range("D33").Value() = 192600
$ws.Range("D34").Value() = 0
$ws.Range("D35").Value() = 192600
$ws.Range("D38").Value() = 43463
$ws.Range("D41").Value() = 143100
$ws.Range("D42").Value() = 0
$ws.Range("D43").Value() = 361200
$ws.Range("D44").Value() = 317600
$ws.Range("D45").Value() = 45800
$ws.Range("D46").Value() = 867700
$ws.Range("D47").Value() = 0
$ws.Range("D48").Value() = 130900
$ws.Range("D49").Value() = 1100800
$ws.Range("D50").Value() = 0
$ws.Range("D51").Value() = 0
$ws.Range("D52").Value() = 83700
$ws.Range("D53").Value() = 0
$ws.Range("D54").Value() = 2183100
$ws.Range("D57").Value() = 202300
$ws.Range("D58").Value() = 132500
$ws.Range("D59").Value() = 138300
$ws.Range("D60").Value() = 473100
$ws.Range("D61").Value() = 438000
$ws.Range("D62").Value() = 280400
$ws.Range("D63").Value() = 0
$ws.Range("D64").Value() = 0
$ws.Range("D65").Value() = 0
$ws.Range("D66").Value() = 1197100
$ws.Range("D68").Value() = 0
$ws.Range("D69").Value() = 0
$ws.Range("D70").Value() = 0
$ws.Range("D71").Value() = 0
$ws.Range("D72").Value() = 1169700
$ws.Range("D73").Value() = 0
$ws.Range("D74").Value() = 0
$ws.Range("D75").Value() = 0
$ws.Range("D76").Value() = 986000
$ws.Range("D77").Value() = 0
$ws.Range("D80").Value() = 43463
$ws.Range("D81").Value() = 192600
$ws.Range("D83").Value() = 31500
$ws.Range("D84").Value() = 0
$ws.Range("D85").Value() = 0
$ws.Range("D86").Value() = 0
$ws.Range("D87").Value() = 0
$ws.Range("D88").Value() = 0
$ws.Range("D89").Value() = 97500
$ws.Range("D91").Value() = -21700
$ws.Range("D92").Value() = 0
$ws.Range("D93").Value() = 0
$ws.Range("D94").Value() = -22200
$ws.Range("D96").Value() = -28600
$ws.Range("D97").Value() = 0
$ws.Range("D98").Value() = 0
$ws.Range("D99").Value() = 0
$ws.Range("D100").Value() = -404500
$ws.Range("D101").Value() = -8700
$ws.Range("D102").Value() = -337900

# A few rows received revised (restated) figures for the two most recent prior periods as well
$ws.Range("E17").Value() = 2318400
$ws.Range("F17").Value() = 2348900
$ws.Range("E18").Value() = 31600
$ws.Range("F18").Value() = 145700
$ws.Range("E20").Value() = -42200
$ws.Range("F20").Value() = -35200
$ws.Range("E32").Value() = 42200
$ws.Range("F32").Value() = 35200
